$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, new Price (D), new Volume(1h) (E), new Hora (G).
# Empty string means "no change for this cell" (matches the source diff).
$updates = @(
    @{Row=2; D="305.62"; E="1.22%"; G="20"},
    @{Row=3; D="35.93"; E="0.76%"; G="20"},
    @{Row=4; D="5.017"; E="-1.17%"; G="20"},
    @{Row=5; D="0.08076"; E="0.36%"; G="20"},
    @{Row=6; D="1.947"; E="0.33%"; G="20"},
    @{Row=7; D="4.143"; E="2.38%"; G="20"},
    @{Row=8; D="7.845"; E="0.87%"; G="20"},
    @{Row=9; D="0.9309"; E="0.52%"; G="20"},
    @{Row=10; D="0.1256"; E="-20.19%"; G="20"},
    @{Row=11; D="0.1914"; E="0.32%"; G="20"},
    @{Row=12; D="0.09195"; E="2.63%"; G="20"},
    @{Row=13; D="0.03512"; E=""; G="20"},
    @{Row=14; D="0.09946"; E="0.58%"; G="20"},
    @{Row=15; D="0.001416"; E="1.42%"; G="20"},
    @{Row=16; D="0.006695"; E="16.74%"; G="20"},
    @{Row=17; D="3.616"; E="2.27%"; G="20"},
    @{Row=18; D="3.006"; E="4.36%"; G="20"},
    @{Row=19; D="0.3436"; E="-0.21%"; G="20"},
    @{Row=20; D="5.172"; E="1.95%"; G="20"},
    @{Row=21; D=""; E="-0.41%"; G="20"},
    @{Row=22; D="0.2532"; E="5.61%"; G="20"},
    @{Row=23; D="0.04406"; E="-2.00%"; G="20"},
    @{Row=24; D="0.001236"; E="2.20%"; G="20"},
    @{Row=25; D="0.004719"; E="-1.26%"; G="20"},
    @{Row=26; D=""; E="5.87%"; G="20"},
    @{Row=27; D="0.0003131"; E="3.70%"; G="20"},
    @{Row=28; D=""; E=""; G="20"},
    @{Row=29; D=""; E=""; G="20"},
    @{Row=30; D=""; E=""; G="20"},
    @{Row=31; D=""; E=""; G="20"},
    @{Row=32; D=""; E=""; G="20"},
    @{Row=33; D=""; E=""; G="20"},
    @{Row=34; D=""; E=""; G="20"},
    @{Row=35; D=""; E=""; G="20"},
    @{Row=36; D=""; E=""; G="20"},
    @{Row=37; D=""; E=""; G="20"},
    @{Row=38; D=""; E=""; G="20"},
    @{Row=39; D="0.01963"; E="6.08%"; G="20"},
    @{Row=40; D="0.05182"; E="7.99%"; G="20"},
    @{Row=41; D="0.007589"; E="4.06%"; G="20"},
    @{Row=42; D="0.01013"; E="-4.52%"; G="20"},
    @{Row=43; D="0.1372"; E="2.81%"; G="20"},
    @{Row=44; D="0.002102"; E="-0.30%"; G="20"},
    @{Row=45; D=""; E="9.94%"; G="20"},
    @{Row=46; D="0.00006395"; E="2.68%"; G="20"},
    @{Row=47; D=""; E="0.23%"; G="20"},
    @{Row=48; D="65.22"; E="0.86%"; G="20"},
    @{Row=49; D="0.001602"; E="-3.48%"; G="20"},
    @{Row=50; D="0.00002102"; E="0.23%"; G="20"},
    @{Row=51; D="0.0002002"; E="0.23%"; G="20"}
)

foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("D", "E", "G")) {
        $newVal = $u[$col]
        if ([string]::IsNullOrEmpty($newVal)) { continue }
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value = $newVal
        $cell.Style = "Normal"
    }
}
